$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns B..N between row 13 and row 14
# (column A, the numeric index, stays on its original row)
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N")

foreach ($col in $cols) {
    $addr13 = "$col" + "13"
    $addr14 = "$col" + "14"

    $val13 = $ws.Range($addr13).Value2
    $val14 = $ws.Range($addr14).Value2

    $ws.Range($addr13).Value = $val14
    $ws.Range($addr14).Value = $val13
}
